$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows (3-21) down to (4-22)
$ws.Rows("3").Insert()

# Populate the newly inserted row 3 with the new data record
$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "Terminal La Palmera de La Serena"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44631
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100101
$ws.Range("H3").Value = "Berries"
$ws.Range("I3").Value = 100101001
$ws.Range("J3").Value = "Arándano (blue)"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 160
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 6500
$ws.Range("P3").Value = 6250
$ws.Range("Q3").Value = "`$/bandeja 2 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 3125
$ws.Range("T3").Value = 2
